$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply the bordered / centered style (same as column A data cells) to all of B2:C31 ---
$ws.Range("A3").Copy()
$ws.Range("B2:C31").PasteSpecial(-4122)   # xlPasteFormats

# --- 2. Enter text values in the exact order the author originally typed them,
#        so the shared-string table indices line up with the target workbook ---
$ws.Range("B3").Value  = "yes"
$ws.Range("C2").Value  = "-"
$ws.Range("C4").Value  = "-"
$ws.Range("C6").Value  = "-"
$ws.Range("C7").Value  = "-"
$ws.Range("C8").Value  = "-"
$ws.Range("C10").Value = "-"
$ws.Range("C5").Value  = "DD/MM/YYYY"
$ws.Range("C9").Value  = "dS/m"
$ws.Range("C11").Value = "ppm"
$ws.Range("C12").Value = "ppm"
$ws.Range("C13").Value = "ppm"
$ws.Range("C14").Value = "ppm"
$ws.Range("C15").Value = "ppm"
$ws.Range("C16").Value = "ppm"
$ws.Range("C17").Value = "ppm"
$ws.Range("C18").Value = "ppm"
$ws.Range("C19").Value = "ppm"
$ws.Range("C20").Value = "ppm"
$ws.Range("C21").Value = "ppm"
$ws.Range("C22").Value = "ppm"
$ws.Range("C23").Value = "ppm"
$ws.Range("C24").Value = "ppm"
$ws.Range("C25").Value = "ppm"
$ws.Range("C26").Value = "ppm"
$ws.Range("C27").Value = "ppm"
$ws.Range("C28").Value = "ppm"
$ws.Range("C29").Value = "ppm"
$ws.Range("C30").Value = "ppm"
$ws.Range("C31").Value = "CFU/ml"
$ws.Range("B7").Value  = "Irrigation"
$ws.Range("B6").Value  = "previous"

# --- 3. Numeric values ---
$ws.Range("B2").Value  = 100
$ws.Range("B4").Value  = 15
$ws.Range("B8").Value  = 2
$ws.Range("B9").Value  = 0.3
$ws.Range("B10").Value = 7.7
$ws.Range("B11").Value = 1.8
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 3
$ws.Range("B16").Value = 48
$ws.Range("B17").Value = 8.2
$ws.Range("B18").Value = 28
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("B25").Value = 0
$ws.Range("B26").Value = 0
$ws.Range("B27").Value = 0
$ws.Range("B28").Value = 0
$ws.Range("B29").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("B31").Value = 0

# --- 4. Date cell B5: numeric date serial, formatted with the short-date built-in format (14) ---
$ws.Range("B5").Value = 43235
$ws.Range("B5").NumberFormat = "mm-dd-yy"

# --- 5. B2 gets its own distinct style (border+center, General format explicitly re-applied) ---
$ws.Range("B2").NumberFormat = "General"

# --- 6. Data validation list on B7 ---
$ws.Range("B7").Validation.Add(3, 1, 1, "#REF!") | Out-Null

# --- 7. Selection state saved in the sheet view ---
$ws.Range("B3").Select() | Out-Null
